$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.154748678207397
$ws.Range("B1").Value = 2.00465989112854
$ws.Range("C1").Value = 5.553575038909912
$ws.Range("D1").Value = 0.739995002746582
$ws.Range("E1").Value = 0.8718262910842896
